$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 807, shifting existing rows 807:848 down to 808:849
$ws.Rows.Item(807).Insert()

# Populate the newly inserted row 807 with the new data point.
# Format column A as text first so the date-looking string "2026/02/11" is
# stored as a literal string rather than being auto-converted to a date
# serial number, then clear the formatting again so the cell ends up with
# the same (default / no explicit style) appearance as its neighbours.
$ws.Range("A807").NumberFormat = "@"
$ws.Range("A807").Value = "2026/02/11"
$ws.Range("A807").ClearFormats()

$ws.Range("B807").Value = "水"
$ws.Range("C807").Value = 17
$ws.Range("D807").Value = 201
